# Add the "naca0015" symmetric airfoil profile sheet (NACA 0015 winglet profile)
# as the last sheet in the workbook, mirroring the other imported-profile sheets
# (Initial MH-61, Sheet2, centre_profile.txt, profile_35mm_away, MH-45 reflex 3,
# MH45 rootprofile).

$wb = $excel.ActiveWorkbook

$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)

# Insert the new sheet right after the current last sheet ("MH45 rootprofile"),
# so it becomes the new last / active tab.
$naca = $wb.Worksheets.Add($null, $lastSheet)
$naca.Name = "naca0015"

# NACA 0015 symmetric profile coordinates (x, y, z) as imported from
# .../payload/Airfoil/database/naca0015.txt
$data = @(
  @(30.0, 0.0474, 0),
  @(28.5, 0.3024, 0),
  @(27.0, 0.543, 0),
  @(24.0, 0.9837, 0),
  @(21.0, 1.374, 0),
  @(18.0, 1.7112, 0),
  @(15.0, 1.9851, 0),
  @(12.0, 2.1761999999999997, 0),
  @(9.0, 2.2506, 0),
  @(7.5, 2.2281, 0),
  @(6.0, 2.1516, 0),
  @(4.5, 2.0046, 0),
  @(3.0, 1.7559, 0),
  @(2.25, 1.575, 0),
  @(1.5, 1.3329, 0),
  @(0.75, 0.9804, 0),
  @(0.375, 0.7101, 0),
  @(0, 0, 0),
  @(0.375, -0.7101, 0),
  @(0.75, -0.9804, 0),
  @(1.5, -1.3329, 0),
  @(2.25, -1.575, 0),
  @(3.0, -1.7559, 0),
  @(4.5, -2.0046, 0),
  @(6.0, -2.1516, 0),
  @(7.5, -2.2281, 0),
  @(9.0, -2.2506, 0),
  @(12.0, -2.1761999999999997, 0),
  @(15.0, -1.9851, 0),
  @(18.0, -1.7112, 0),
  @(21.0, -1.374, 0),
  @(24.0, -0.9837, 0),
  @(27.0, -0.543, 0),
  @(28.5, -0.3024, 0),
  @(30.0, -0.0474, 0)

)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $naca.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Approximate the "best fit" column widths Excel's text-import would have applied.
$naca.Columns.Item(1).ColumnWidth = 7.1640625
$naca.Columns.Item(2).ColumnWidth = 8.6640625

# Sheet-scoped defined name pointing at the imported range, like the workbook's
# other profile sheets.
$naca.Names.Add("naca0015", "=naca0015!`$A`$1:`$B`$35")

# Match the author's final selection/active-tab state.
$naca.Range("E10").Select()
